$wb = $excel.ActiveWorkbook

# --- weibull ---
$ws = $wb.Worksheets.Item("weibull")
$ws.Range("B2").Value = -2.79738339481437
$ws.Range("C2").Value = 0.171698196289065
$ws.Range("B3").Value = -0.0139801110920569
$ws.Range("C3").Value = 0.105889677952872

# --- lognormal ---
$ws = $wb.Worksheets.Item("lognormal")
$ws.Range("B2").Value = 2.20329872621066
$ws.Range("C2").Value = 0.251472192731928
$ws.Range("B3").Value = -0.949589611386702
$ws.Range("C3").Value = 0.120483913853294

# --- llogis ---
$ws = $wb.Worksheets.Item("llogis")
$ws.Range("B2").Value = -2.29189352078826
$ws.Range("C2").Value = 0.134661724775328
$ws.Range("B3").Value = 1.63959728057625
$ws.Range("C3").Value = 0.243901739439504

# --- gompertz ---
$ws = $wb.Worksheets.Item("gompertz")
$ws.Range("B2").Value = -2.49902757553759
$ws.Range("C2").Value = 0.151796300191924
$ws.Range("B3").Value = -0.0260432363531547
$ws.Range("C3").Value = 0.012122266923568

# --- weibull cov ---
$ws = $wb.Worksheets.Item("weibull cov")
$ws.Range("A2").Value = 0.0294802706089184
$ws.Range("B2").Value = -0.0116199924086907
$ws.Range("A3").Value = -0.0116199924086907
$ws.Range("B3").Value = 0.011212623896963

# --- lognormal cov ---
$ws = $wb.Worksheets.Item("lognormal cov")
$ws.Range("A2").Value = 0.0632382637174039
$ws.Range("B2").Value = -0.0269335727041458
$ws.Range("A3").Value = -0.0269335727041458
$ws.Range("B3").Value = 0.0145163734974079

# --- llogis cov ---
$ws = $wb.Worksheets.Item("llogis cov")
$ws.Range("A2").Value = 0.0181337801194661
$ws.Range("B2").Value = 0.0106255584561674
$ws.Range("A3").Value = 0.0106255584561674
$ws.Range("B3").Value = 0.0594880585016157

# --- gompertz cov ---
$ws = $wb.Worksheets.Item("gompertz cov")
$ws.Range("A2").Value = 0.0230421167519567
$ws.Range("B2").Value = -0.000973927502218067
$ws.Range("A3").Value = -0.000973927502218067
$ws.Range("B3").Value = 0.000146949355366231
